$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source cells store plain text (dates as "YYYY/MM/DD" strings and
# numbers as formatted text) rather than real Excel dates/numbers.
# Assigning a date- or number-looking string via .Value would make Excel
# auto-convert it to a serial date / number, so we temporarily force the
# cell to Text format, write the literal string, then clear the
# formatting override again so the cell keeps the workbook's original
# (unstyled) look while retaining the text value.
function Set-TextValue {
    param($cell, [string]$text)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# Row -> (new date, new EBITDA text-or-$null-if-unchanged)
$updates = @{
    2  = @{ A = "2025/12/02"; B = "5.10" }
    8  = @{ A = "2025/12/02"; B = "7.71" }
    14 = @{ A = "2025/12/02"; B = $null }
    20 = @{ A = "2025/12/02"; B = "12.51" }
    26 = @{ A = "2025/12/02"; B = "10.01" }
    32 = @{ A = "2025/12/02"; B = "26.27" }
    38 = @{ A = "2025/12/02"; B = $null }
    44 = @{ A = "2025/12/02"; B = $null }
    50 = @{ A = "2025/12/02"; B = "11.62" }
    56 = @{ A = "2025/12/02"; B = "33.33" }
    62 = @{ A = "2025/12/02"; B = "11.42" }
    68 = @{ A = "2025/12/02"; B = "11.88" }
    74 = @{ A = "2025/12/02"; B = "16.14" }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]

    $aCell = $ws.Cells.Item($row, 1)
    Set-TextValue $aCell $vals.A

    if ($null -ne $vals.B) {
        $bCell = $ws.Cells.Item($row, 2)
        Set-TextValue $bCell $vals.B
    }
}
